$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 175 (existing rows 175-220 shift down to 178-223)
$ws.Range("A175:A177").EntireRow.Insert()

# Boilerplate columns shared by every data row in this sheet
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112028
$categoria = "Sandia"
$variedad  = "Sin especificar"
$kgOUnid   = 1
$clasif    = "Hortaliza"

$newRows = @(
    @{ Row=175; Fecha=44943; Calidad="Extra";    Volumen=500; PMin=3000; PMax=3000; PProm=3000; Unidad="$/unidad"; Origen="Región del Maule"; PKg=3000 },
    @{ Row=176; Fecha=44943; Calidad="Primera";  Volumen=500; PMin=2500; PMax=2500; PProm=2500; Unidad="$/unidad"; Origen="Región del Maule"; PKg=2500 },
    @{ Row=177; Fecha=44943; Calidad="Segunda";  Volumen=500; PMin=2000; PMax=2000; PProm=2000; Unidad="$/unidad"; Origen="Región del Maule"; PKg=2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgOUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
